$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 19
$prevRow = $newRow - 1

# Values for the new data row (matches the pattern of the preceding rows).
$dateText  = "2025/11/28"
$gameName  = "逃离鸭科夫"
$modCount  = 1275

# Leading apostrophe forces the date-looking string to be stored as literal
# text (so it round-trips as "2025/11/28" instead of being auto-converted
# into a date serial number).
$ws.Cells.Item($newRow, 1).Value = "'" + $dateText
$ws.Cells.Item($newRow, 2).Value = $gameName
$ws.Cells.Item($newRow, 3).Value = $modCount

# Copy the formatting (centered alignment style used by the data rows)
# from the row above onto the newly added row.
$ws.Range("A" + $prevRow + ":C" + $prevRow).Copy()
$ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
